$wb = $excel.ActiveWorkbook

# ALC row 9: Distill, My Heart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 11366.889
$ws.Range("I9").Value = 15217.167
$ws.Range("J9").Value = 3666.3333
$ws.Range("K9").Value = 15217.167
$ws.Range("L9").Value = 3666.3333
$ws.Range("M9").Value = -15048.167
$ws.Range("N9").Value = -4004.3333

# ALC row 32: Automata for the People
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7527.143
$ws.Range("I32").Value = 6000
$ws.Range("J32").Value = 8138
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 8138
$ws.Range("M32").Value = -5674
$ws.Range("N32").Value = -8790

# ALC row 80: Cleansing the Wicked Humours
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 796.7143
$ws.Range("I80").Value = 158
$ws.Range("J80").Value = 970.9091
$ws.Range("K80").Value = 474
$ws.Range("L80").Value = 2912.7273
$ws.Range("M80").Value = 524
$ws.Range("N80").Value = -4908.7273

# ALC row 83: Washing Away the Sins (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 796.7143
$ws.Range("I83").Value = 158
$ws.Range("J83").Value = 970.9091
$ws.Range("K83").Value = 1422
$ws.Range("L83").Value = 8738.1819
$ws.Range("M83").Value = 3570
$ws.Range("N83").Value = -18722.1819

# ALC row 111: An Eye for Healing
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 20876.375
$ws.Range("I111").Value = 20900
$ws.Range("J111").Value = 20868.5
$ws.Range("K111").Value = 62700
$ws.Range("L111").Value = 62605.5
$ws.Range("M111").Value = -59633
$ws.Range("N111").Value = -68739.5

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1850.6875
$ws.Range("I137").Value = 1909.5714
$ws.Range("J137").Value = 1804.8889
$ws.Range("K137").Value = 5728.7142
$ws.Range("L137").Value = 5414.6667
$ws.Range("M137").Value = -3178.7142
$ws.Range("N137").Value = -10514.6667

# ARM row 61: Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3838.982
$ws.Range("I61").Value = 2096.4285
$ws.Range("J61").Value = 5646.074
$ws.Range("K61").Value = 2096.4285
$ws.Range("L61").Value = 5646.074
$ws.Range("M61").Value = -1884.4285
$ws.Range("N61").Value = -6070.074

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10998.143
$ws.Range("I74").Value = 1745.3334
$ws.Range("J74").Value = 34130.168
$ws.Range("K74").Value = 1745.3334
$ws.Range("L74").Value = 34130.168
$ws.Range("M74").Value = -871.3334
$ws.Range("N74").Value = -35878.168

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 10998.143
$ws.Range("I77").Value = 1745.3334
$ws.Range("J77").Value = 34130.168
$ws.Range("K77").Value = 8726.666999999999
$ws.Range("L77").Value = 170650.84
$ws.Range("M77").Value = -4358.666999999999
$ws.Range("N77").Value = -179386.84

# ARM row 136: Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3838.982
$ws.Range("I136").Value = 2096.4285
$ws.Range("J136").Value = 5646.074
$ws.Range("K136").Value = 6289.2855
$ws.Range("L136").Value = 16938.222
$ws.Range("M136").Value = -3739.2855
$ws.Range("N136").Value = -22038.222

# BSM row 20: Smelt and Dealt
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4500.8335
$ws.Range("I20").Value = 3507.0667
$ws.Range("J20").Value = 5494.6
$ws.Range("K20").Value = 3507.0667
$ws.Range("L20").Value = 5494.6
$ws.Range("M20").Value = -3260.0667
$ws.Range("N20").Value = -5988.6

# CRP row 80: The Long Armillae of the Law
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 31775
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 31775
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 31775
$ws.Range("N80").Value = -34021

# CRP row 83: Wooden Ambitions (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 31775
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 31775
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 95325
$ws.Range("N83").Value = -106557

# CRP row 86: Birch, Please
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7277.1665
$ws.Range("I86").Value = 6556
$ws.Range("J86").Value = 7998.3335
$ws.Range("K86").Value = 6556
$ws.Range("L86").Value = 7998.3335
$ws.Range("M86").Value = -5433
$ws.Range("N86").Value = -10244.3335

# CRP row 89: Built This City on Blocks and Soul (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7277.1665
$ws.Range("I89").Value = 6556
$ws.Range("J89").Value = 7998.3335
$ws.Range("K89").Value = 32780
$ws.Range("L89").Value = 39991.6675
$ws.Range("M89").Value = -27164
$ws.Range("N89").Value = -51223.6675

# CUL row 23: Sweet Smell of Success
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 399.76923
$ws.Range("I23").Value = 382.75
$ws.Range("J23").Value = 407.33334
$ws.Range("K23").Value = 1148.25
$ws.Range("L23").Value = 1222.00002
$ws.Range("M23").Value = -913.25
$ws.Range("N23").Value = -1692.00002

# CUL row 113: Can't Eat Just One
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 621.4761999999999
$ws.Range("I113").Value = 675.8570999999999
$ws.Range("J113").Value = 594.2857
$ws.Range("K113").Value = 2027.5713
$ws.Range("L113").Value = 1782.8571
$ws.Range("M113").Value = 142.4287000000002
$ws.Range("N113").Value = -6122.8571

# CUL row 129: Comfort Food
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 9001621
$ws.Range("I129").Value = 24751094
$ws.Range("J129").Value = 1922.7142
$ws.Range("K129").Value = 74253282
$ws.Range("L129").Value = 5768.142599999999
$ws.Range("M129").Value = -74248282
$ws.Range("N129").Value = -15768.1426

# GSM row 104: Speak Softly and Carry a Metal Rod
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 30000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 30000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

# GSM row 107: Whetstones for the Workers
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 965.75
$ws.Range("I107").Value = 821.3333
$ws.Range("J107").Value = 1399
$ws.Range("K107").Value = 821.3333
$ws.Range("L107").Value = 1399
$ws.Range("M107").Value = 1098.6667
$ws.Range("N107").Value = -5239

# LTW row 16: Saddle Sore
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2559.0476
$ws.Range("I16").Value = 2670.9473
$ws.Range("J16").Value = 1496
$ws.Range("K16").Value = 2670.9473
$ws.Range("L16").Value = 1496
$ws.Range("M16").Value = -2500.9473
$ws.Range("N16").Value = -1836

# LTW row 80: Don't Sweat the Small Fry
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 30062.5
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 40125
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 40125
$ws.Range("M80").Value = -18877
$ws.Range("N80").Value = -42371

# LTW row 83: It's All in the Wrists (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H83").Value = 30062.5
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 40125
$ws.Range("K83").Value = 60000
$ws.Range("L83").Value = 120375
$ws.Range("M83").Value = -54384
$ws.Range("N83").Value = -131607

# LTW row 121: A Shoe In
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 99997.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 99997.5
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 99997.5
$ws.Range("N121").Value = -103491.5

# LTW row 132: Tenets of Tanning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5177.846
$ws.Range("I132").Value = 4980.2
$ws.Range("J132").Value = 5836.6665
$ws.Range("K132").Value = 14940.6
$ws.Range("L132").Value = 17509.9995
$ws.Range("M132").Value = -12410.6
$ws.Range("N132").Value = -22569.9995

# WVR row 81: Where the Dragonflies, the Net Catches
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9373.579
$ws.Range("I81").Value = 21859.8
$ws.Range("J81").Value = 4914.2144
$ws.Range("K81").Value = 43719.6
$ws.Range("L81").Value = 9828.4288
$ws.Range("M81").Value = -42658.6
$ws.Range("N81").Value = -11950.4288

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 9373.579
$ws.Range("I84").Value = 21859.8
$ws.Range("J84").Value = 4914.2144
$ws.Range("K84").Value = 218598
$ws.Range("L84").Value = 49142.144
$ws.Range("M84").Value = -213294
$ws.Range("N84").Value = -59750.144

# WVR row 88: The Hat List
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 16057
$ws.Range("I88").Value = 17085.5
$ws.Range("J88").Value = 14000
$ws.Range("K88").Value = 17085.5
$ws.Range("L88").Value = 14000
$ws.Range("M88").Value = -16679.5
$ws.Range("N88").Value = -14812

# WVR row 91: Knight Incognito (L)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H91").Value = 16057
$ws.Range("I91").Value = 17085.5
$ws.Range("J91").Value = 14000
$ws.Range("K91").Value = 17085.5
$ws.Range("L91").Value = 14000
$ws.Range("M91").Value = -15681.5
$ws.Range("N91").Value = -16808

# WVR row 107: Flax Wax
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 84351.414
$ws.Range("I107").Value = 1117.4286
$ws.Range("J107").Value = 200879
$ws.Range("K107").Value = 3352.2858
$ws.Range("L107").Value = 602637
$ws.Range("M107").Value = -1432.2858
$ws.Range("N107").Value = -606477

# WVR row 113: A Tender Table
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2012.9584
$ws.Range("I113").Value = 1286.909
$ws.Range("J113").Value = 9999.5
$ws.Range("K113").Value = 3860.727
$ws.Range("L113").Value = 29998.5
$ws.Range("M113").Value = -1690.727
$ws.Range("N113").Value = -34340
